# Auto-generated edit script applying the Fruta/Hortaliza weekly data refresh.
# Re-assigns Fecha/Volumen/Precio values per row and appends one new record (row 44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("P2").Value = 972
# Row 3
$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("P3").Value = 972
# Row 4
$ws.Range("D4").Value = (Get-Date -Year 2022 -Month 4 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861
# Row 5
$ws.Range("D5").Value = (Get-Date -Year 2022 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15500
$ws.Range("P5").Value = 861
# Row 6
$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I6").Value = "Primera"
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("P6").Value = 944
# Row 7
$ws.Range("D7").Value = (Get-Date -Year 2022 -Month 3 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 861
# Row 8
$ws.Range("D8").Value = (Get-Date -Year 2022 -Month 7 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("P8").Value = 972
# Row 9
$ws.Range("D9").Value = (Get-Date -Year 2022 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 861
# Row 10
$ws.Range("D10").Value = (Get-Date -Year 2022 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J10").Value = 120
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 17000
$ws.Range("P10").Value = 944
# Row 11
$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 972
# Row 12
$ws.Range("D12").Value = (Get-Date -Year 2021 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I12").Value = "Segunda"
$ws.Range("K12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("P12").Value = 833
# Row 13
$ws.Range("D13").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 861
# Row 14
$ws.Range("D14").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J14").Value = 100
# Row 15
$ws.Range("D15").Value = (Get-Date -Year 2022 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 16500
$ws.Range("P15").Value = 917
# Row 16
$ws.Range("D16").Value = (Get-Date -Year 2023 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15500
$ws.Range("P16").Value = 861
# Row 17
$ws.Range("D17").Value = (Get-Date -Year 2022 -Month 4 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J17").Value = 80
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 861
# Row 18
$ws.Range("D18").Value = (Get-Date -Year 2022 -Month 8 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J18").Value = 100
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 972
# Row 19
$ws.Range("D19").Value = (Get-Date -Year 2022 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("P19").Value = 972
# Row 20
$ws.Range("D20").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J20").Value = 100
# Row 21
$ws.Range("D21").Value = (Get-Date -Year 2022 -Month 8 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("P21").Value = 972
# Row 22
$ws.Range("D22").Value = (Get-Date -Year 2022 -Month 7 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("P22").Value = 972
# Row 23
$ws.Range("D23").Value = (Get-Date -Year 2023 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 13500
$ws.Range("P23").Value = 750
# Row 24
$ws.Range("D24").Value = (Get-Date -Year 2022 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J24").Value = 60
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("P24").Value = 833
# Row 25
$ws.Range("D25").Value = (Get-Date -Year 2022 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 16000
$ws.Range("M25").Value = 16500
$ws.Range("P25").Value = 917
# Row 26
$ws.Range("D26").Value = (Get-Date -Year 2022 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17500
$ws.Range("P26").Value = 972
# Row 27
$ws.Range("D27").Value = (Get-Date -Year 2022 -Month 9 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K27").Value = 15000
$ws.Range("M27").Value = 15000
$ws.Range("P27").Value = 833
# Row 29
$ws.Range("D29").Value = (Get-Date -Year 2022 -Month 8 -Day 18 -Hour 0 -Minute 0 -Second 0)
# Row 30
$ws.Range("D30").Value = (Get-Date -Year 2022 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J30").Value = 60
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("P30").Value = 861
# Row 31
$ws.Range("D31").Value = (Get-Date -Year 2023 -Month 4 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("P31").Value = 833
# Row 32
$ws.Range("D32").Value = (Get-Date -Year 2022 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = 15500
$ws.Range("P32").Value = 861
# Row 33
$ws.Range("D33").Value = (Get-Date -Year 2022 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 806
# Row 34
$ws.Range("D34").Value = (Get-Date -Year 2022 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J34").Value = 60
# Row 35
$ws.Range("D35").Value = (Get-Date -Year 2022 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J35").Value = 60
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 14500
$ws.Range("P35").Value = 806
# Row 36
$ws.Range("D36").Value = (Get-Date -Year 2023 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 14000
$ws.Range("M36").Value = 14000
$ws.Range("P36").Value = 778
# Row 37
$ws.Range("D37").Value = (Get-Date -Year 2022 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K37").Value = 14000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 14500
$ws.Range("P37").Value = 806
# Row 38
$ws.Range("D38").Value = (Get-Date -Year 2022 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K38").Value = 17000
$ws.Range("L38").Value = 17000
$ws.Range("M38").Value = 17000
$ws.Range("P38").Value = 944
# Row 39
$ws.Range("D39").Value = (Get-Date -Year 2023 -Month 4 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J39").Value = 60
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 15000
$ws.Range("P39").Value = 833
# Row 40
$ws.Range("D40").Value = (Get-Date -Year 2023 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K40").Value = 16000
$ws.Range("L40").Value = 16000
$ws.Range("M40").Value = 16000
$ws.Range("P40").Value = 889
# Row 41
$ws.Range("D41").Value = (Get-Date -Year 2022 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J41").Value = 80
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 15500
$ws.Range("P41").Value = 861
# Row 42
$ws.Range("D42").Value = (Get-Date -Year 2022 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J42").Value = 100
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = 15500
$ws.Range("P42").Value = 861
# Row 43
$ws.Range("D43").Value = (Get-Date -Year 2022 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)
# Row 44
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = (Get-Date -Year 2023 -Month 4 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino dulce"
$ws.Range("H44").Value = "Cultivar IV Región"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 80
$ws.Range("K44").Value = 15000
$ws.Range("L44").Value = 16000
$ws.Range("M44").Value = 15500
$ws.Range("N44").Value = "`$/bandeja 18 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 861
$ws.Range("Q44").Value = 18
$ws.Range("R44").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D44").NumberFormat = $ws.Range("D43").NumberFormat
